$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.125615"
$ws.Range("H2").Value = [double]"0.25123"
$ws.Range("I2").Value = [double]"0.02647478672532295"
$ws.Range("J2").Value = [double]"0.01780700335556722"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.5"
$ws.Range("M2").Value = [double]"0.071358"
$ws.Range("N2").Value = [double]"0.142716"
$ws.Range("O2").Value = [double]"0.01919591193090569"
$ws.Range("P2").Value = [double]"0.01411929935366186"
$ws.Range("Q2").Value = [double]"0.008963635170000001"
$ws.Range("R2").Value = [double]"0.03585454068"
$ws.Range("S2").Value = [double]"0.0005082076743688106"
$ws.Range("T2").Value = [double]"0.0002514224109689147"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.125615"
$ws.Range("H3").Value = [double]"0.25123"
$ws.Range("I3").Value = [double]"0.02647478672532295"
$ws.Range("J3").Value = [double]"0.01780700335556722"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"2.655685"
$ws.Range("N3").Value = [double]"7.967055"
$ws.Range("O3").Value = [double]"0.7144019644080171"
$ws.Range("P3").Value = [double]"0.7882033865305114"
$ws.Range("Q3").Value = [double]"0.333593871275"
$ws.Range("R3").Value = [double]"2.00156322765"
$ws.Range("S3").Value = [double]"0.01891363964385401"
$ws.Range("T3").Value = [double]"0.01403554034881826"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.125615"
$ws.Range("H4").Value = [double]"0.25123"
$ws.Range("I4").Value = [double]"0.02647478672532295"
$ws.Range("J4").Value = [double]"0.01780700335556722"
$ws.Range("K4").Value = [double]"2"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"0.9728370000000001"
$ws.Range("N4").Value = [double]"1.945674"
$ws.Range("O4").Value = [double]"0.2617014683024538"
$ws.Range("P4").Value = [double]"0.1924910567184946"
$ws.Range("Q4").Value = [double]"0.122202919755"
$ws.Range("R4").Value = [double]"0.4888116790200001"
$ws.Range("S4").Value = [double]"0.006928490559011331"
$ws.Range("T4").Value = [double]"0.003427688892902914"

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = [double]"2"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"0.125615"
$ws.Range("H5").Value = [double]"0.25123"
$ws.Range("I5").Value = [double]"0.02647478672532295"
$ws.Range("J5").Value = [double]"0.01780700335556722"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.017474"
$ws.Range("N5").Value = [double]"0.052422"
$ws.Range("O5").Value = [double]"0.004700655358623364"
$ws.Range("P5").Value = [double]"0.005186257397332197"
$ws.Range("Q5").Value = [double]"0.00219499651"
$ws.Range("R5").Value = [double]"0.01316997906"
$ws.Range("S5").Value = [double]"0.0001244488480888"
$ws.Range("T5").Value = [double]"9.235170287712975E-05"

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Ror2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"4.619088000000001"
$ws.Range("H6").Value = [double]"13.857264"
$ws.Range("I6").Value = [double]"0.9735252132746771"
$ws.Range("J6").Value = [double]"0.9821929966444328"
$ws.Range("K6").Value = [double]"1"
$ws.Range("L6").Value = [double]"0.5"
$ws.Range("M6").Value = [double]"0.071358"
$ws.Range("N6").Value = [double]"0.142716"
$ws.Range("O6").Value = [double]"0.01919591193090569"
$ws.Range("P6").Value = [double]"0.01411929935366186"
$ws.Range("Q6").Value = [double]"0.3296088815040001"
$ws.Range("R6").Value = [double]"1.977653289024"
$ws.Range("S6").Value = [double]"0.01868770425653688"
$ws.Range("T6").Value = [double]"0.01386787694269294"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Ror2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"4.619088000000001"
$ws.Range("H7").Value = [double]"13.857264"
$ws.Range("I7").Value = [double]"0.9735252132746771"
$ws.Range("J7").Value = [double]"0.9821929966444328"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"2.655685"
$ws.Range("N7").Value = [double]"7.967055"
$ws.Range("O7").Value = [double]"0.7144019644080171"
$ws.Range("P7").Value = [double]"0.7882033865305114"
$ws.Range("Q7").Value = [double]"12.26684271528"
$ws.Range("R7").Value = [double]"110.40158443752"
$ws.Range("S7").Value = [double]"0.6954883247641631"
$ws.Range("T7").Value = [double]"0.7741678461816931"

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Ror2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"4.619088000000001"
$ws.Range("H8").Value = [double]"13.857264"
$ws.Range("I8").Value = [double]"0.9735252132746771"
$ws.Range("J8").Value = [double]"0.9821929966444328"
$ws.Range("K8").Value = [double]"2"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"0.9728370000000001"
$ws.Range("N8").Value = [double]"1.945674"
$ws.Range("O8").Value = [double]"0.2617014683024538"
$ws.Range("P8").Value = [double]"0.1924910567184946"
$ws.Range("Q8").Value = [double]"4.493619712656001"
$ws.Range("R8").Value = [double]"26.961718275936"
$ws.Range("S8").Value = [double]"0.2547729777434425"
$ws.Range("T8").Value = [double]"0.1890633678255917"

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Ror2"
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"4.619088000000001"
$ws.Range("H9").Value = [double]"13.857264"
$ws.Range("I9").Value = [double]"0.9735252132746771"
$ws.Range("J9").Value = [double]"0.9821929966444328"
$ws.Range("K9").Value = [double]"1"
$ws.Range("L9").Value = [double]"0.3333333333333333"
$ws.Range("M9").Value = [double]"0.017474"
$ws.Range("N9").Value = [double]"0.052422"
$ws.Range("O9").Value = [double]"0.004700655358623364"
$ws.Range("P9").Value = [double]"0.005186257397332197"
$ws.Range("Q9").Value = [double]"0.08071394371200001"
$ws.Range("R9").Value = [double]"0.7264254934080001"
$ws.Range("S9").Value = [double]"0.004576206510534564"
$ws.Range("T9").Value = [double]"0.005093905694455067"
